$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 3).Value = 766310
$ws.Cells.Item(2, 4).Value = 155968
$ws.Cells.Item(2, 5).Value = 1429187466

$ws.Cells.Item(10, 3).Value = 345533
$ws.Cells.Item(10, 4).Value = 64164
$ws.Cells.Item(10, 5).Value = 1817711330

$ws.Cells.Item(13, 3).Value = 187834
$ws.Cells.Item(13, 4).Value = 33262
$ws.Cells.Item(13, 5).Value = 1165224910

$ws.Cells.Item(54, 3).Value = 75191
$ws.Cells.Item(54, 4).Value = 14914
$ws.Cells.Item(54, 5).Value = 361048094

$ws.Cells.Item(78, 3).Value = 178440
$ws.Cells.Item(78, 4).Value = 34685
$ws.Cells.Item(78, 5).Value = 892473543

$ws.Cells.Item(81, 3).Value = 88350
$ws.Cells.Item(81, 4).Value = 16598
$ws.Cells.Item(81, 5).Value = 499614738

$ws.Cells.Item(88, 3).Value = 71264
$ws.Cells.Item(88, 4).Value = 12436
$ws.Cells.Item(88, 5).Value = 110294557

$ws.Cells.Item(91, 3).Value = 18848
$ws.Cells.Item(91, 4).Value = 3383
$ws.Cells.Item(91, 5).Value = 75117834

$ws.Cells.Item(93, 3).Value = 16923
$ws.Cells.Item(93, 4).Value = 2945
$ws.Cells.Item(93, 5).Value = 50452265

$ws.Cells.Item(121, 3).Value = 1306126
$ws.Cells.Item(121, 4).Value = 220385
$ws.Cells.Item(121, 5).Value = 2274557857

$ws.Cells.Item(122, 3).Value = 364
$ws.Cells.Item(122, 4).Value = 48
$ws.Cells.Item(122, 5).Value = 1156270

$ws.Cells.Item(129, 3).Value = 633338
$ws.Cells.Item(129, 4).Value = 104966
$ws.Cells.Item(129, 5).Value = 3426710411

$ws.Cells.Item(132, 3).Value = 585622
$ws.Cells.Item(132, 4).Value = 90778
$ws.Cells.Item(132, 5).Value = 3461185017

$ws.Cells.Item(136, 3).Value = 26675
$ws.Cells.Item(136, 4).Value = 4272
$ws.Cells.Item(136, 5).Value = 143540532

$ws.Cells.Item(139, 3).Value = 76639
$ws.Cells.Item(139, 4).Value = 17494
$ws.Cells.Item(139, 5).Value = 114131597

$ws.Cells.Item(144, 3).Value = 25068
$ws.Cells.Item(144, 4).Value = 6170
$ws.Cells.Item(144, 5).Value = 92362703

$ws.Cells.Item(151, 3).Value = 39921
$ws.Cells.Item(151, 4).Value = 7155
$ws.Cells.Item(151, 5).Value = 60359605

$ws.Cells.Item(154, 3).Value = 18438
$ws.Cells.Item(154, 4).Value = 3296
$ws.Cells.Item(154, 5).Value = 72659539

$ws.Cells.Item(156, 3).Value = 12397
$ws.Cells.Item(156, 4).Value = 2144
$ws.Cells.Item(156, 5).Value = 40042868

$ws.Cells.Item(158, 3).Value = 717
$ws.Cells.Item(158, 4).Value = 132
$ws.Cells.Item(158, 5).Value = 1762994

$ws.Cells.Item(159, 3).Value = 43847
$ws.Cells.Item(159, 4).Value = 5415
$ws.Cells.Item(159, 5).Value = 101312981

$ws.Cells.Item(178, 3).Value = 515876
$ws.Cells.Item(178, 4).Value = 115380
$ws.Cells.Item(178, 5).Value = 891189200

$ws.Cells.Item(207, 3).Value = 154660
$ws.Cells.Item(207, 4).Value = 27118
$ws.Cells.Item(207, 5).Value = 753642756
